$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "28.141.48"
$ws.Range("E2").Value = "  -0.10%  "

# Row 3
$ws.Range("D3").Value = "1.760.66"
$ws.Range("E3").Value = "  -2.65%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "334.40"
$ws.Range("E5").Value = "  -1.51%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9982"
$ws.Range("E6").Value = "  -0.38%  "

# Row 7
$ws.Range("E7").Value = "  -3.68%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3364"
$ws.Range("E8").Value = "  -3.45%  "

# Row 9
$ws.Range("E9").Value = "  -5.66%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.130"
$ws.Range("E10").Value = "  -4.88%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07215"
$ws.Range("E11").Value = "  -4.36%  "

# Row 12
$ws.Range("B12").Value = "Solana"
$ws.Range("C12").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.73"
$ws.Range("E12").Value = "  +2.97%  "

# Row 13
$ws.Range("B13").Value = "BinanceUSD"
$ws.Range("C13").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.0000"
$ws.Range("E13").Value = "  -0.17%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.181"
$ws.Range("E14").Value = "  -5.01%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.217"
$ws.Range("E15").Value = "  +0.98%  "

# Row 16
$ws.Range("D16").Value = "1.758.59"
$ws.Range("E16").Value = "  -2.93%  "

# Row 17
$ws.Range("E17").Value = "  -4.36%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06576"
$ws.Range("E18").Value = "  -1.82%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "80.93"
$ws.Range("E19").Value = "  -4.60%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9984"
$ws.Range("E20").Value = "  -0.25%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.01"
$ws.Range("E21").Value = "  -4.12%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.272"
$ws.Range("E22").Value = "  -4.65%  "

# Row 23
$ws.Range("D23").Value = "28.132.55"
$ws.Range("E23").Value = "  -0.08%  "

# Row 24
$ws.Range("E24").Value = "  -6.32%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.396"
$ws.Range("E25").Value = "  -0.66%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "153.00"
$ws.Range("E26").Value = "  -0.37%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.90"

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.332"
$ws.Range("E28").Value = "  -7.85%  "

# Row 29
$ws.Range("D29").Value = "1.959.89"
$ws.Range("E29").Value = "  -2.86%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.265"
$ws.Range("E30").Value = "  -15.05%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "132.23"
$ws.Range("E31").Value = "  -2.32%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.016"
$ws.Range("E32").Value = "  -0.21%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.821"
$ws.Range("E33").Value = "  -5.53%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08770"
$ws.Range("E34").Value = "  -0.88%  "

# Row 35
$ws.Range("E35").Value = "  -5.70%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02346"
$ws.Range("E36").Value = "  -3.01%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6658"
$ws.Range("E37").Value = "  -4.01%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06218"
$ws.Range("E38").Value = "  -5.09%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.168"
$ws.Range("E39").Value = "  -5.17%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2118"
$ws.Range("E40").Value = "  -4.17%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.221"
$ws.Range("E41").Value = "  -2.85%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.454"
$ws.Range("E42").Value = "  -9.64%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.010"

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9980"
$ws.Range("E44").Value = "  -0.16%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.72"
$ws.Range("E45").Value = "  -6.14%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6067"
$ws.Range("E46").Value = "  -5.46%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.812"
$ws.Range("E47").Value = "  -1.55%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "129.72"
$ws.Range("E48").Value = "  -1.25%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.014"
$ws.Range("E49").Value = "  -5.98%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.196"
$ws.Range("E50").Value = "  +3.39%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07213"
$ws.Range("E51").Value = "  +0.26%  "
